$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# Five pairs of fixtures that were recorded on the same day got their rows
# swapped (re-sorted) by the source-data refresh. Column A (the running id)
# keeps counting 0,1,2,... in row order, so it is left untouched; every other
# column (B..AC) of the two rows trades places. The date in column E is the
# same for both rows of a pair, so nothing changes there either.
# ----------------------------------------------------------------------

# Swap rows 26 and 27
$ws.Range("B26").Value2 = 6732773
$ws.Range("F26").Value2 = 'Suduva Marijampole'
$ws.Range("G26").Value2 = 'Hegelmann Litauen'
$ws.Range("H26").Value2 = 0
$ws.Range("I26").Value2 = 1
$ws.Range("J26").Value2 = 'A'
$ws.Range("K26").Value2 = 5
$ws.Range("L26").Value2 = 3.8
$ws.Range("M26").Value2 = 1.533
$ws.Range("N26").Value2 = 5
$ws.Range("O26").Value2 = 4.2
$ws.Range("P26").Value2 = 1.533
$ws.Range("Q26").Value2 = 1
$ws.Range("R26").Value2 = 1.875
$ws.Range("S26").Value2 = 1.925
$ws.Range("T26").Value2 = 2.5
$ws.Range("U26").Value2 = 1.9
$ws.Range("V26").Value2 = 1.9
$ws.Range("W26").Value2 = -1
$ws.Range("X26").Value2 = -1
$ws.Range("Y26").Value2 = 0.5329999999999999
$ws.Range("Z26").Value2 = 0
$ws.Range("AA26").Value2 = -0
$ws.Range("AB26").Value2 = -1
$ws.Range("AC26").Value2 = 0.8999999999999999
$ws.Range("B27").Value2 = 6732711
$ws.Range("F27").Value2 = 'Banga Gargzdai'
$ws.Range("G27").Value2 = 'FK Zalgiris Vilnius'
$ws.Range("H27").Value2 = 1
$ws.Range("I27").Value2 = 4
$ws.Range("J27").Value2 = 'A'
$ws.Range("K27").Value2 = 5
$ws.Range("L27").Value2 = 3.6
$ws.Range("M27").Value2 = 1.571
$ws.Range("N27").Value2 = 11
$ws.Range("O27").Value2 = 4.75
$ws.Range("P27").Value2 = 1.25
$ws.Range("Q27").Value2 = 1.5
$ws.Range("R27").Value2 = 1.975
$ws.Range("S27").Value2 = 1.825
$ws.Range("T27").Value2 = 2.5
$ws.Range("U27").Value2 = 1.8
$ws.Range("V27").Value2 = 2
$ws.Range("W27").Value2 = -1
$ws.Range("X27").Value2 = -1
$ws.Range("Y27").Value2 = 0.25
$ws.Range("Z27").Value2 = -1
$ws.Range("AA27").Value2 = 0.825
$ws.Range("AB27").Value2 = 0.8
$ws.Range("AC27").Value2 = -1

# Swap rows 89 and 90
$ws.Range("B89").Value2 = 6732827
$ws.Range("F89").Value2 = 'FK Dziugas Telsiai'
$ws.Range("G89").Value2 = 'FK Kauno Zalgiris'
$ws.Range("H89").Value2 = 0
$ws.Range("I89").Value2 = 2
$ws.Range("J89").Value2 = 'A'
$ws.Range("K89").Value2 = 6
$ws.Range("L89").Value2 = 3.9
$ws.Range("M89").Value2 = 1.444
$ws.Range("N89").Value2 = 4.75
$ws.Range("O89").Value2 = 3.6
$ws.Range("P89").Value2 = 1.65
$ws.Range("Q89").Value2 = 0.75
$ws.Range("R89").Value2 = 1.9
$ws.Range("S89").Value2 = 1.9
$ws.Range("T89").Value2 = 2.5
$ws.Range("U89").Value2 = 1.95
$ws.Range("V89").Value2 = 1.85
$ws.Range("W89").Value2 = -1
$ws.Range("X89").Value2 = -1
$ws.Range("Y89").Value2 = 0.6499999999999999
$ws.Range("Z89").Value2 = -1
$ws.Range("AA89").Value2 = 0.8999999999999999
$ws.Range("AB89").Value2 = -1
$ws.Range("AC89").Value2 = 0.8500000000000001
$ws.Range("B90").Value2 = 7326568
$ws.Range("F90").Value2 = 'Hegelmann Litauen'
$ws.Range("G90").Value2 = 'Panevezys'
$ws.Range("H90").Value2 = 0
$ws.Range("I90").Value2 = 0
$ws.Range("J90").Value2 = 'D'
$ws.Range("K90").Value2 = 2.375
$ws.Range("L90").Value2 = 3.2
$ws.Range("M90").Value2 = 2.625
$ws.Range("N90").Value2 = 2.7
$ws.Range("O90").Value2 = 3.2
$ws.Range("P90").Value2 = 2.3
$ws.Range("Q90").Value2 = 0
$ws.Range("R90").Value2 = 2.05
$ws.Range("S90").Value2 = 1.75
$ws.Range("T90").Value2 = 2.25
$ws.Range("U90").Value2 = 1.875
$ws.Range("V90").Value2 = 1.925
$ws.Range("W90").Value2 = -1
$ws.Range("X90").Value2 = 2.2
$ws.Range("Y90").Value2 = -1
$ws.Range("Z90").Value2 = 0
$ws.Range("AA90").Value2 = -0
$ws.Range("AB90").Value2 = -1
$ws.Range("AC90").Value2 = 0.925

# Swap rows 100 and 104
$ws.Range("B100").Value2 = 6732836
$ws.Range("F100").Value2 = 'FK Siauliai'
$ws.Range("G100").Value2 = 'Banga Gargzdai'
$ws.Range("H100").Value2 = 3
$ws.Range("I100").Value2 = 0
$ws.Range("J100").Value2 = 'H'
$ws.Range("K100").Value2 = 1.222
$ws.Range("L100").Value2 = 5.5
$ws.Range("M100").Value2 = 9
$ws.Range("N100").Value2 = 1.363
$ws.Range("O100").Value2 = 4.5
$ws.Range("P100").Value2 = 7
$ws.Range("Q100").Value2 = -1.25
$ws.Range("R100").Value2 = 1.9
$ws.Range("S100").Value2 = 1.9
$ws.Range("T100").Value2 = 2.5
$ws.Range("U100").Value2 = 1.975
$ws.Range("V100").Value2 = 1.825
$ws.Range("W100").Value2 = 0.363
$ws.Range("X100").Value2 = -1
$ws.Range("Y100").Value2 = -1
$ws.Range("Z100").Value2 = 0.8999999999999999
$ws.Range("AA100").Value2 = -1
$ws.Range("AB100").Value2 = 0.9750000000000001
$ws.Range("AC100").Value2 = -1
$ws.Range("B104").Value2 = 6732837
$ws.Range("F104").Value2 = 'Suduva Marijampole'
$ws.Range("G104").Value2 = 'FK Riteriai'
$ws.Range("H104").Value2 = 0
$ws.Range("I104").Value2 = 3
$ws.Range("J104").Value2 = 'A'
$ws.Range("K104").Value2 = 3.6
$ws.Range("L104").Value2 = 3.6
$ws.Range("M104").Value2 = 1.8
$ws.Range("N104").Value2 = 3
$ws.Range("O104").Value2 = 3.6
$ws.Range("P104").Value2 = 2
$ws.Range("Q104").Value2 = 0.25
$ws.Range("R104").Value2 = 2
$ws.Range("S104").Value2 = 1.8
$ws.Range("T104").Value2 = 2.5
$ws.Range("U104").Value2 = 1.975
$ws.Range("V104").Value2 = 1.825
$ws.Range("W104").Value2 = -1
$ws.Range("X104").Value2 = -1
$ws.Range("Y104").Value2 = 1
$ws.Range("Z104").Value2 = -1
$ws.Range("AA104").Value2 = 0.8
$ws.Range("AB104").Value2 = 0.9750000000000001
$ws.Range("AC104").Value2 = -1

# Swap rows 101 and 103
$ws.Range("B101").Value2 = 6732834
$ws.Range("F101").Value2 = 'Panevezys'
$ws.Range("G101").Value2 = 'FK Dziugas Telsiai'
$ws.Range("H101").Value2 = 0
$ws.Range("I101").Value2 = 0
$ws.Range("J101").Value2 = 'D'
$ws.Range("K101").Value2 = 1.25
$ws.Range("L101").Value2 = 5.5
$ws.Range("M101").Value2 = 7.5
$ws.Range("N101").Value2 = 1.45
$ws.Range("O101").Value2 = 4.5
$ws.Range("P101").Value2 = 5
$ws.Range("Q101").Value2 = -1
$ws.Range("R101").Value2 = 1.775
$ws.Range("S101").Value2 = 2.025
$ws.Range("T101").Value2 = 2.5
$ws.Range("U101").Value2 = 1.875
$ws.Range("V101").Value2 = 1.925
$ws.Range("W101").Value2 = -1
$ws.Range("X101").Value2 = 3.5
$ws.Range("Y101").Value2 = -1
$ws.Range("Z101").Value2 = -1
$ws.Range("AA101").Value2 = 1.025
$ws.Range("AB101").Value2 = -1
$ws.Range("AC101").Value2 = 0.925
$ws.Range("B103").Value2 = 7465686
$ws.Range("F103").Value2 = 'FK Kauno Zalgiris'
$ws.Range("G103").Value2 = 'Hegelmann Litauen'
$ws.Range("H103").Value2 = 4
$ws.Range("I103").Value2 = 2
$ws.Range("J103").Value2 = 'H'
$ws.Range("K103").Value2 = 2.3
$ws.Range("L103").Value2 = 4
$ws.Range("M103").Value2 = 2.3
$ws.Range("N103").Value2 = 2.55
$ws.Range("O103").Value2 = 4
$ws.Range("P103").Value2 = 2.2
$ws.Range("Q103").Value2 = 0.25
$ws.Range("R103").Value2 = 1.8
$ws.Range("S103").Value2 = 2
$ws.Range("T103").Value2 = 2.75
$ws.Range("U103").Value2 = 1.85
$ws.Range("V103").Value2 = 1.95
$ws.Range("W103").Value2 = 1.55
$ws.Range("X103").Value2 = -1
$ws.Range("Y103").Value2 = -1
$ws.Range("Z103").Value2 = 0.8
$ws.Range("AA103").Value2 = -1
$ws.Range("AB103").Value2 = 0.8500000000000001
$ws.Range("AC103").Value2 = -1

# Swap rows 117 and 118
$ws.Range("B117").Value2 = 7862911
$ws.Range("F117").Value2 = 'Hegelmann Litauen'
$ws.Range("G117").Value2 = 'FK Siauliai'
$ws.Range("H117").Value2 = 2
$ws.Range("I117").Value2 = 2
$ws.Range("J117").Value2 = 'D'
$ws.Range("K117").Value2 = 2.15
$ws.Range("L117").Value2 = 3.1
$ws.Range("M117").Value2 = 3.1
$ws.Range("N117").Value2 = 2.45
$ws.Range("O117").Value2 = 2.9
$ws.Range("P117").Value2 = 3
$ws.Range("Q117").Value2 = 0
$ws.Range("R117").Value2 = 1.725
$ws.Range("S117").Value2 = 2.075
$ws.Range("T117").Value2 = 2.5
$ws.Range("U117").Value2 = 2.025
$ws.Range("V117").Value2 = 1.775
$ws.Range("W117").Value2 = -1
$ws.Range("X117").Value2 = 1.9
$ws.Range("Y117").Value2 = -1
$ws.Range("Z117").Value2 = 0
$ws.Range("AA117").Value2 = -0
$ws.Range("AB117").Value2 = 1.025
$ws.Range("AC117").Value2 = -1
$ws.Range("B118").Value2 = 7862036
$ws.Range("F118").Value2 = 'Banga Gargzdai'
$ws.Range("G118").Value2 = 'FK Zalgiris Vilnius'
$ws.Range("H118").Value2 = 1
$ws.Range("I118").Value2 = 4
$ws.Range("J118").Value2 = 'A'
$ws.Range("K118").Value2 = 8
$ws.Range("L118").Value2 = 4.5
$ws.Range("M118").Value2 = 1.3
$ws.Range("N118").Value2 = 6.5
$ws.Range("O118").Value2 = 4.5
$ws.Range("P118").Value2 = 1.333
$ws.Range("Q118").Value2 = 1.25
$ws.Range("R118").Value2 = 2
$ws.Range("S118").Value2 = 1.8
$ws.Range("T118").Value2 = 2.5
$ws.Range("U118").Value2 = 1.825
$ws.Range("V118").Value2 = 1.975
$ws.Range("W118").Value2 = -1
$ws.Range("X118").Value2 = -1
$ws.Range("Y118").Value2 = 0.333
$ws.Range("Z118").Value2 = -1
$ws.Range("AA118").Value2 = 0.8
$ws.Range("AB118").Value2 = 0.825
$ws.Range("AC118").Value2 = -1

# ----------------------------------------------------------------------
# Three brand-new fixtures were added at the bottom of the table, rows 125-127
# (dimension grows from A1:AC124 to A1:AC127).
# ----------------------------------------------------------------------
# Row 125
$ws.Range("A125").Value2 = 123
$ws.Range("B125").Value2 = 7862915
$ws.Range("C125").Value2 = 'Lithuania A Lyga'
$ws.Range("D125").Value2 = 'Lithuania A Lyga'
$ws.Range("E125").Value2 = 45380.58333333334
$ws.Range("F125").Value2 = 'Banga Gargzdai'
$ws.Range("G125").Value2 = 'FK Dziugas Telsiai'
$ws.Range("H125").Value2 = 0
$ws.Range("I125").Value2 = 2
$ws.Range("J125").Value2 = 'A'
$ws.Range("K125").Value2 = 2.6
$ws.Range("L125").Value2 = 2.9
$ws.Range("M125").Value2 = 2.625
$ws.Range("N125").Value2 = 2.6
$ws.Range("O125").Value2 = 2.8
$ws.Range("P125").Value2 = 2.75
$ws.Range("Q125").Value2 = 0
$ws.Range("R125").Value2 = 1.825
$ws.Range("S125").Value2 = 1.975
$ws.Range("T125").Value2 = 2
$ws.Range("U125").Value2 = 1.975
$ws.Range("V125").Value2 = 1.825
$ws.Range("W125").Value2 = -1
$ws.Range("X125").Value2 = -1
$ws.Range("Y125").Value2 = 1.75
$ws.Range("Z125").Value2 = -1
$ws.Range("AA125").Value2 = 0.9750000000000001
$ws.Range("AB125").Value2 = 0
$ws.Range("AC125").Value2 = 0
$ws.Range("A2").Copy()
$ws.Range("A125").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E125").PasteSpecial(-4122)

# Row 126
$ws.Range("A126").Value2 = 124
$ws.Range("B126").Value2 = 7862916
$ws.Range("C126").Value2 = 'Lithuania A Lyga'
$ws.Range("D126").Value2 = 'Lithuania A Lyga'
$ws.Range("E126").Value2 = 45381.41666666666
$ws.Range("F126").Value2 = 'Hegelmann Litauen'
$ws.Range("G126").Value2 = 'Panevezys'
$ws.Range("H126").Value2 = 4
$ws.Range("I126").Value2 = 2
$ws.Range("J126").Value2 = 'H'
$ws.Range("K126").Value2 = 2.5
$ws.Range("L126").Value2 = 3
$ws.Range("M126").Value2 = 2.6
$ws.Range("N126").Value2 = 2.6
$ws.Range("O126").Value2 = 2.9
$ws.Range("P126").Value2 = 2.6
$ws.Range("Q126").Value2 = 0
$ws.Range("R126").Value2 = 1.9
$ws.Range("S126").Value2 = 1.9
$ws.Range("T126").Value2 = 2.25
$ws.Range("U126").Value2 = 1.925
$ws.Range("V126").Value2 = 1.875
$ws.Range("W126").Value2 = 1.6
$ws.Range("X126").Value2 = -1
$ws.Range("Y126").Value2 = -1
$ws.Range("Z126").Value2 = 0.8999999999999999
$ws.Range("AA126").Value2 = -1
$ws.Range("AB126").Value2 = 0.925
$ws.Range("AC126").Value2 = -1
$ws.Range("A2").Copy()
$ws.Range("A126").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E126").PasteSpecial(-4122)

# Row 127
$ws.Range("A127").Value2 = 125
$ws.Range("B127").Value2 = 7862917
$ws.Range("C127").Value2 = 'Lithuania A Lyga'
$ws.Range("D127").Value2 = 'Lithuania A Lyga'
$ws.Range("E127").Value2 = 45381.5
$ws.Range("F127").Value2 = 'FK Siauliai'
$ws.Range("G127").Value2 = 'FK Transinvest'
$ws.Range("H127").Value2 = 0
$ws.Range("I127").Value2 = 1
$ws.Range("J127").Value2 = 'A'
$ws.Range("K127").Value2 = 1.533
$ws.Range("L127").Value2 = 3.75
$ws.Range("M127").Value2 = 5
$ws.Range("N127").Value2 = 1.65
$ws.Range("O127").Value2 = 3.4
$ws.Range("P127").Value2 = 4.5
$ws.Range("Q127").Value2 = -0.75
$ws.Range("R127").Value2 = 1.875
$ws.Range("S127").Value2 = 1.925
$ws.Range("T127").Value2 = 2.75
$ws.Range("U127").Value2 = 1.95
$ws.Range("V127").Value2 = 1.85
$ws.Range("W127").Value2 = -1
$ws.Range("X127").Value2 = -1
$ws.Range("Y127").Value2 = 3.5
$ws.Range("Z127").Value2 = -1
$ws.Range("AA127").Value2 = 0.925
$ws.Range("AB127").Value2 = -1
$ws.Range("AC127").Value2 = 0.8500000000000001
$ws.Range("A2").Copy()
$ws.Range("A127").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E127").PasteSpecial(-4122)

$excel.CutCopyMode = 0
